# Insert a new Title-styled paragraph " 100 node references" as the very
# first paragraph of the document body (before the existing paragraph that
# holds the inline picture).

$d = $word.ActiveDocument

# Collapsed range at the very start of the document's main story.
$r = $d.Range(0, 0)

# Build the new paragraph as a WordprocessingML fragment and insert it
# in one shot, so the paragraph mark/style/run/text all land exactly as
# in the target markup (no extra rsid bookkeeping from piecemeal edits).
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:pPr><w:pStyle w:val="Title"/></w:pPr>' +
            '<w:r><w:t xml:space="preserve"> 100 node references</w:t></w:r>' +
            '</w:p>'

[void]$r.InsertXML($titleXml)
